$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('A2').Value = 'Barrendero/A - Turno Nocturno o Vespertino'
$ws.Range('B2').Value = 'Idea Market Solutions'
$ws.Range('C2').Value = 'Iztapalapa, Ciudad de México DF'
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '$ 5,437'
$ws.Range('E2').Value = 'Buscamos tu talento como Barrendero/a para turno Nocturno o Vespertino:Zona de trabajo: Central de abastos IztapalapaBarrenderos Nocturno (Horario 5 pm a 5 am)Barrenderos Vespertino (Horario de 10 am a 7:00 pm)Lunes a domingo descanso entre semana*Contratacion inmediata*Contar con documentacion (Acta de nacimiento, Comprobante de clabe interbancaria es decir cuenta de banco, nss, curp, ine, constancia de situacion fiscal y comprobante de domicilo)Agenda hoy mismo comunicate al 5 5 796 9 59 09'
$ws.Range('F2').Value = 'https://mx.computrabajo.com/ofertas-de-trabajo/oferta-de-trabajo-de-barrenderoa-turno-nocturno-o-vespertino-en-iztapalapa-8EFACEA56D1028F261373E686DCF3405#lc=ListOffers-Score4-0'

$ws.Range('A3').Value = 'Barrendero Central de Abastos - Iztapalapa'
$ws.Range('B3').Value = 'Idea Market Solutions'
$ws.Range('C3').Value = 'Iztapalapa, Ciudad de México DF'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '$ 10,000.00'
$ws.Range('E3').Value = 'Importante mercado está en busca de tu talentoSOLICITA BARRENDERO MATUTINO O NOCTURNOOFRECEMOS:Sueldo de $8,0000 a $10,000 neto mensualUniforme gratuitoPrestaciones de leyREQUISITOS:Experiencia no necesariaPrimaria concluidaDocumentos completos en copiaLunes a Domingo de 5pm a 5am o 10am a 7pm, 1 descanso rolado a la semanaZona de trabajo: Central de abastos Iztapalapa, CDMX5511 4122 82Si estas interesada/o postúlate con tu CV actualizado'
$ws.Range('F3').Value = 'https://mx.computrabajo.com/ofertas-de-trabajo/oferta-de-trabajo-de-barrendero-central-de-abastos-iztapalapa-en-iztapalapa-D4CB1C1EE38BEE9A61373E686DCF3405#lc=ListOffers-Score4-1'

$ws.Range('A4').Value = 'Barrendero con o sin experiencia / Nocturno y Vespertino - Zona Central de abastos Iztapalapa'
$ws.Range('B4').Value = 'Importante empresa del sector'
$ws.Range('C4').Value = 'Iztapalapa, Ciudad de México DF'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '$9,000'
$ws.Range('E4').Value = 'Se busca barrendero para turnos nocturnos y vespertinos.Buscamos a una persona responsable, con experiencia en limpieza y mantenimiento de espacios comunes.El trabajo requiere desplazamientos cortos y uso de herramientas básicas de limpieza.Se valorará experiencia previa en limpieza, aunque no es imprescindible.Se ofrece contrato temporal con posibilidad de renovación.Se requiere disponibilidad para trabajar en horarios nocturnos y vespertinos.Se busca un perfil con actitud positiva y compromiso con el trabajo bien hecho.Se requiere puntualidad y capacidad para trabajar en equipo.Se ofrece ambiente de trabajo agradable y posibilidad de formación continua.Si tienes interés en mantener espacios limpios y cuidando los detalles, esta oportunidad puede ser para ti.Envía tu currículum con tus datos y experiencia previa en limpieza.Buscas un trabajo estable y con oportunidades de crecimiento profesional.Únete a nuestro equipo y contribuye a mantener nuestros entornos limpios y ordenados.'
$ws.Range('F4').Value = 'https://mx.computrabajo.com/ofertas-de-trabajo/oferta-de-trabajo-de-barrendero-con-o-sin-experiencia-nocturno-y-vespertino-zona-central-de-abastos-iztapalapa-en-iztapalapa-9BB6F3E7844AC9D961373E686DCF3405#lc=ListOffers-Score4-2'

$ws.Range('A5').Value = 'EG Barrendero Vespertino Central de Abastos'
$ws.Range('B5').Value = 'Idea Market Solutions'
$ws.Range('C5').Value = 'Iztapalapa, Ciudad de México DF'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '$ 8,000.00'
$ws.Range('E5').Value = 'Importante empresa de MANTENIMIENTO está en busca de tu talentoSolicita: BARRENDERO VESPERTINOOFRECEMOS:• Sueldo de hasta $8,000• Prestaciones de ley• PAGO SEMANALREQUISITOS:• SIN EXPERIENCIA• Escolaridad mínima• L a D• Horario fijoZONA DE TRABAJO: CENTRAL DE ABASTOS, IZTAPALAPA, CDMXSi estas interesada/o postúlate con tu CV actualizado'
$ws.Range('F5').Value = 'https://mx.computrabajo.com/ofertas-de-trabajo/oferta-de-trabajo-de-eg-barrendero-vespertino-central-de-abastos-en-iztapalapa-B6999002976CCE4B61373E686DCF3405#lc=ListOffers-Score4-3'

$ws.Range('A6').Value = 'EG Barrendero Nocturno Central de Abastos'
$ws.Range('B6').Value = 'Idea Market Solutions'
$ws.Range('C6').Value = 'Iztapalapa, Ciudad de México DF'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '$ 10,000.00'
$ws.Range('E6').Value = 'Importante empresa de MANTENIMIENTO está en busca de tu talentoSolicita: BARRENDERO NOCTURNOOFRECEMOS:• Sueldo de hasta $10,000• Prestaciones de ley• PAGO SEMANALREQUISITOS:• SIN EXPERIENCIA• Escolaridad mínima• L a D• Horario fijoZONA DE TRABAJO: CENTRAL DE ABASTOS, IZTAPALAPA, CDMXSi estas interesada/o postúlate con tu CV actualizado'
$ws.Range('F6').Value = 'https://mx.computrabajo.com/ofertas-de-trabajo/oferta-de-trabajo-de-eg-barrendero-nocturno-central-de-abastos-en-iztapalapa-BFADA104A2B4F0B661373E686DCF3405#lc=ListOffers-Score4-4'

$ws.Range('A7').Value = 'Barrenderos matutino y nocturno con o sin experiencia'
$ws.Range('B7').Value = 'Idea Market Solutions'
$ws.Range('C7').Value = 'Iztapalapa, Ciudad de México DF'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '$ 8,000.00'
$ws.Range('E7').Value = 'Unete al mejor equipo de trabajo como BARRENDER@TURNOS DISPONIBLESmatutino: 9:00 am - 6:00 pmnocturno: 6pm - 6 amLunes a Sábado a Domingo descanso rolado entre semanaSalario: $2,000 semanales ( $2,500 NOCTURNO)+ prestaciones+uniformesSexo indistintoNo necesitas experienciaContratación inmediataZona de trabajo: Central de Abastos CDMX ( IZTAPALAPA)'
$ws.Range('F7').Value = 'https://mx.computrabajo.com/ofertas-de-trabajo/oferta-de-trabajo-de-barrenderos-matutino-y-nocturno-con-o-sin-experiencia-en-iztapalapa-11A7D3B0A60709DB61373E686DCF3405#lc=ListOffers-Score4-5'

$ws.Range('A8').Value = 'EG Barrendero Nocturno Central de Abastos'
$ws.Range('B8').Value = 'Idea Market Solutions'
$ws.Range('C8').Value = 'Iztapalapa, Ciudad de México DF'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '$ 10,000.00'
$ws.Range('E8').Value = 'Importante empresa de MANTENIMIENTO está en busca de tu talentoSolicita: BARRENDERO NOCTURNOOFRECEMOS:• Sueldo de hasta $10,000• Prestaciones de ley• PAGO SEMANALREQUISITOS:• SIN EXPERIENCIA• Escolaridad mínima• L a D• Horario fijoZONA DE TRABAJO: CENTRAL DE ABASTOS, IZTAPALAPA, CDMXSi estas interesada/o postúlate con tu CV actualizado'
$ws.Range('F8').Value = 'https://mx.computrabajo.com/ofertas-de-trabajo/oferta-de-trabajo-de-eg-barrendero-nocturno-central-de-abastos-en-iztapalapa-C9417E1A9ACBCE5A61373E686DCF3405#lc=ListOffers-Score4-6'

$ws.Range('A9').Value = 'EG Barrendero Vespertino Central de Abastos'
$ws.Range('B9').Value = 'Idea Market Solutions'
$ws.Range('C9').Value = 'Iztapalapa, Ciudad de México DF'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '$ 8,000.00'
$ws.Range('E9').Value = 'Importante empresa de MANTENIMIENTO está en busca de tu talentoSolicita: BARRENDERO VESPERTINOOFRECEMOS:• Sueldo de hasta $8,000• Prestaciones de ley• PAGO SEMANALREQUISITOS:• SIN EXPERIENCIA• Escolaridad mínima• L a D• Horario fijoZONA DE TRABAJO: CENTRAL DE ABASTOS, IZTAPALAPA, CDMXSi estas interesada/o postúlate con tu CV actualizado'
$ws.Range('F9').Value = 'https://mx.computrabajo.com/ofertas-de-trabajo/oferta-de-trabajo-de-eg-barrendero-vespertino-central-de-abastos-en-iztapalapa-41FCE64536945C4661373E686DCF3405#lc=ListOffers-Score4-7'

$ws.Range('A10').Value = 'Barrendero - Central de abastos cdmx (iztapalapa)'
$ws.Range('B10').Value = 'Idea Market Solutions'
$ws.Range('C10').Value = 'Iztapalapa, Ciudad de México DF'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '$ 5,437'
$ws.Range('E10').Value = 'Se solicita para Central de Abastos:Contratación inmediata- Barrenderos NocturnoHorario 5 pm a 5 amSueldo: 10 mil pesos mensuales- Barrenderos VespertinoHorario de 10 am a 7:00 pmSueldo: 8000 pesos mensualesBarrido de patios de la centralDocumentos: Acta, ine, curp, nss, constancia sat, comprobante domicilio, estado de cuenta bancario'
$ws.Range('F10').Value = 'https://mx.computrabajo.com/ofertas-de-trabajo/oferta-de-trabajo-de-barrendero-central-de-abastos-cdmx-iztapalapa-en-iztapalapa-001D869701974DF161373E686DCF3405#lc=ListOffers-Score4-8'

$ws.Range('A11').Value = 'Barrendero con turno fijo en la tarde - En Central de Abastos Iztapalapa'
$ws.Range('B11').Value = 'Idea Market Solutions'
$ws.Range('C11').Value = 'Iztapalapa, Ciudad de México DF'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '$ 10,000.00'
$ws.Range('E11').Value = 'Importante empresa dedicada al sector limpieza solicita Barrendero.Requisitos:- Mayor de 18 años- Experiencia no necesariaOfrecemos:- $8,000 mensuales netos (Pago semanal)- Prestaciones de ley- Trabajo de Lunes a Domingo con 1 descanso en la semana- Horario de 10:00 am a 7:00 pmZona de trabajo: Central de Abastos Iztapalapa.Actividades : acompañamiento de camionesLevantar basura de cajones de tolvaLevantar basura que cae de los camiones al momento de la cargaLevantar basura al paso del recorridoInteresados que cuenten con documentos (Acta de nacimiento, Cuenta de Banco, Comprobante de domicilio, INE, CURP, IMSS, RFC) postularse por éste medio, y un ejecutivo te contactará por whatsapp.'
$ws.Range('F11').Value = 'https://mx.computrabajo.com/ofertas-de-trabajo/oferta-de-trabajo-de-barrendero-con-turno-fijo-en-la-tarde-en-central-de-abastos-iztapalapa-en-iztapalapa-33BD941F42C610F461373E686DCF3405#lc=ListOffers-Score4-9'

$ws.Range('A12').Value = 'Barrendero Nocturno - En Central de Abastos Iztapalapa'
$ws.Range('B12').Value = 'Idea Market Solutions'
$ws.Range('C12').Value = 'Iztapalapa, Ciudad de México DF'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '$ 10,000.00'
$ws.Range('E12').Value = 'Importante empresa dedicada al sector limpieza solicita Barrendero para el turno Nocturno.Requisitos:- Mayor de 18 años- Experiencia no necesariaOfrecemos:- $10,000 mensuales netos (Pago semanal)- Prestaciones de ley- Trabajo de Lunes a Domingo con 1 descanso en la semana- Horario de 5:00 pm a 5:00 amZona de trabajo: Central de Abastos Iztapalapa.Actividades : acompañamiento de camionesLevantar basura de cajones de tolvaLevantar basura que cae de los camiones al momento de la cargaLevantar basura al paso del recorridoInteresados que cuenten con documentos (Acta de nacimiento, Cuenta de Banco, Comprobante de domicilio, INE, CURP, IMSS, RFC) postularse por éste medio, y un ejecutivo te contactará por whatsapp.'
$ws.Range('F12').Value = 'https://mx.computrabajo.com/ofertas-de-trabajo/oferta-de-trabajo-de-barrendero-nocturno-en-central-de-abastos-iztapalapa-en-iztapalapa-0F3EFD8ABAA464EE61373E686DCF3405#lc=ListOffers-Score4-10'

$ws.Range('A13').Value = 'Barrenderos / $2,000 a la semana Central de Abastos CDMX'
$ws.Range('B13').Value = 'Idea Market Solutions'
$ws.Range('C13').Value = 'Iztapalapa, Ciudad de México DF'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '$2,000'
$ws.Range('E13').Value = 'SOLICITAMOS POR EXPANSION: BARRENDEROSZONA: IZTAPALAPA (CENTRAL DE ABASTOS)SUELDO $2,000 semanales libresPagos semanalesContarás con Prestaciones de LeyHorario: 9:00 am a 6:00 pmLunes a Domingo, 1 día de descanso entre semana'
$ws.Range('F13').Value = 'https://mx.computrabajo.com/ofertas-de-trabajo/oferta-de-trabajo-de-barrenderos-2000-a-la-semana-central-de-abastos-cdmx-en-iztapalapa-75FC0F398FF27B2C61373E686DCF3405#lc=ListOffers-Score4-11'

$ws.Range('A14').Value = 'Auxiliar limpieza - Central de Abastos'
$ws.Range('B14').Value = 'Idea Market Solutions'
$ws.Range('C14').Value = 'Iztapalapa, Ciudad de México DF'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '$ 8,000.00'
$ws.Range('E14').Value = 'Estamos en búsqueda de tu talento como:TENEMOS DISPONIBLES 2 VACANTES DE AUXILIAR DE LIMPIEZA:BARRENDEROTurno: Nocturno (De 6:00PM a 6:00AM)Sueldo: $10,000 mensualesBARRENDEROTurno: Vespertino (De 10:00AM a 7:00AM)Sueldo: $8,000 mensuales¡NO NECESITAS EXPERIENCIA! - INGRESO INMEDIATOSi estas interesada/o postúlate con tu CV actualizado'
$ws.Range('F14').Value = 'https://mx.computrabajo.com/ofertas-de-trabajo/oferta-de-trabajo-de-auxiliar-limpieza-central-de-abastos-en-iztapalapa-37C5E76C32BA148461373E686DCF3405#lc=ListOffers-Score3-12'

$ws.Range('A15').Value = 'Auxiliar de Limpieza sin Experiencia para Iztapalapa Turno Nocturno - $10,000 mensual Contratación el mismo dia'
$ws.Range('B15').Value = 'Idea Market Solutions'
$ws.Range('C15').Value = 'Iztapalapa, Ciudad de México DF'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '$10,000'
$ws.Range('E15').Value = '"Importante marca de Limpieza está en busca de tu talento"Solicita Axiliar de limpieza (Barrendero) sin ExperienciaRequisitos:Tener disponibilidad de HorarioBuena ActitudLunes a Domingo con descanso entre semanaEdad: entre 18 a 55 añosNo contar con Infonavit de preferenciaSin tatuajes o perforaciones de preferenciaZona de trabajo: IztapalapaOFRECEMOS:Sueldo base Neto semana $2500 semanalNocturno 5pm a 5am $2500 semanalesPago SemanalVacante FijaSi estas interesada/o postúlate con tu CV actualizado'
$ws.Range('F15').Value = 'https://mx.computrabajo.com/ofertas-de-trabajo/oferta-de-trabajo-de-auxiliar-de-limpieza-sin-experiencia-para-iztapalapa-turno-nocturno-10000-mensual-contratacion-el-mismo-dia-en-iztapalapa-F37FE05104C1EFD961373E686DCF3405#lc=ListOffers-Score3-13'

$ws.Range('A16').Value = 'Auxiliar de Limpieza sin Experiencia para Iztapalapa Turno Intermedio - $8000 mensual Contratación el mismo dia'
$ws.Range('B16').Value = 'Idea Market Solutions'
$ws.Range('C16').Value = 'Iztapalapa, Ciudad de México DF'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '$8000'
$ws.Range('E16').Value = '"Importante marca de Limpieza está en busca de tu talento"Solicita Axiliar de limpieza (Barrendero) sin ExperienciaRequisitos:Tener disponibilidad de HorarioBuena ActitudLunes a Domingo con descanso entre semana 10 am a 7pmEdad: entre 18 a 55 añosNo contar con Infonavit de preferenciaSin tatuajes o perforaciones de preferenciaZona de trabajo: IztapalapaOFRECEMOS:Sueldo base Neto semanal $2000 semanalIntermedio 10am a 7pm $2000 semanalesPago SemanalVacante FijaSi estas interesada/o postúlate con tu CV actualizado'
$ws.Range('F16').Value = 'https://mx.computrabajo.com/ofertas-de-trabajo/oferta-de-trabajo-de-auxiliar-de-limpieza-sin-experiencia-para-iztapalapa-turno-intermedio-8000-mensual-contratacion-el-mismo-dia-en-iztapalapa-21DBF84AE98E8B6161373E686DCF3405#lc=ListOffers-Score3-14'

$ws.Range('A17').Value = 'Auxiliares de limpieza contratación inmediata - Iztapalapa Central de abastos'
$ws.Range('B17').Value = 'Idea Market Solutions'
$ws.Range('C17').Value = 'Iztapalapa, Ciudad de México DF'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '$ 10,000.00'
$ws.Range('E17').Value = 'Importante empresa de residuos está en busca de tu talentoSolicitamos personal de limpieza (diversas áreas)BarrenderosAuxiliares de limpieza de pasillosLavadores de camionesZONA IZTAPALAPAOFRECEMOS:Sueldo semanal (dependiendo de la posición)Horario dependiendo de la operaciónPagos semanalesRequisitos:Disponibilidad de horarioExperiencia no necesaria o mínimaZona de trabajoCENTRAL DE ABASTOS IZTAPALAPASi estas interesada/o postúlate con tu CV actualizado o marca al 5580109495'
$ws.Range('F17').Value = 'https://mx.computrabajo.com/ofertas-de-trabajo/oferta-de-trabajo-de-auxiliares-de-limpieza-contratacion-inmediata-iztapalapa-central-de-abastos-en-iztapalapa-65BB4A5C627E4E3761373E686DCF3405#lc=ListOffers-Score3-15'

